$wb = $excel.ActiveWorkbook

# 1) Rename the third sheet "Sheet3" -> "C"
$wsC = $wb.Worksheets.Item(3)
$wsC.Name = "C"

# 2) Work on the first sheet ("大杂烩") - add start time / status columns
$ws1 = $wb.Worksheets.Item(1)

# B2: start time for the book currently being read (row 2 / "Ongoing")
$ws1.Range("B2").Value = 42534
$ws1.Range("B2").NumberFormat = "mm-dd-yy"

# D2: status of the currently-read book
$ws1.Range("D2").Value = "Ongoing"

# D3:D12: status of the rest of the books (not started yet)
$ws1.Range("D3:D12").Value = "Backlog"

# Update the selection on sheet 1 to reflect the newly filled-in status column
$ws1.Range("D3:D12").Select() | Out-Null

# 3) Sheet "C" gets a simple selection update too
$wsC.Range("D34").Select() | Out-Null

# Reselect sheet 1 so it stays the active/visible tab
$ws1.Select() | Out-Null

Write-Host "Done"
